$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Kyle Kuzma", "PF", "Washington Wizards"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
